$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
# NumberFormat/Style dance forces values like "582.88" or "0.330" to be stored
# as text (matching the source feed formatting, incl. trailing zeros) instead of
# being auto-coerced to numbers by Excel, then restores the default "Normal" style
# so no stray number-format style is left on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.034.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.466.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -2.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.467.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.330"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.917.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.918.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("E17").Value = "  -4.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.427.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.47%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -52.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.567.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "517.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0903"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.91%  "
$ws.Range("E34").Value = "  -6.74%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -8.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.62%  "
$ws.Range("E40").Value = "  -6.50%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("E43").Value = "  -7.14%  "
$ws.Range("E44").Value = "  -7.42%  "
$ws.Range("E45").Value = "  -7.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.60%  "
$ws.Range("E48").Value = "  -7.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0254"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -12.56%  "
$ws.Range("E51").Value = "  -7.74%  "
